# Automatische test-sync: 2025-06-19 21:44:50
# Append a new mail-log entry to the "Logs" sheet and bump the matching
# category tally on the "Dashboard" sheet.

$wb = $excel.ActiveWorkbook
$logs = $wb.Worksheets.Item("Logs")
$dashboard = $wb.Worksheets.Item("Dashboard")

$newRow = 27

$logs.Cells.Item($newRow, 1).Value = "Sponsoraanvraag"
$logs.Cells.Item($newRow, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item($newRow, 3).Value = "Zou uw bedrijf bereid zijn om ons sportevenement te sponsoren?"
$logs.Cells.Item($newRow, 4).Value = "Samenwerking / Partnerverzoek"
$logs.Cells.Item($newRow, 6).Value = "2025-06-19 21:44:10"
$logs.Cells.Item($newRow, 7).Value = "Nee"

# Bump the "Samenwerking / Partnerverzoek" count on the Dashboard sheet.
$dashboard.Range("B2").Value = 8

# Extend the conditional-formatting ranges to cover the new row.
$dFc = $logs.Range("D2:D26").FormatConditions
$dFc.Item(1).ModifyAppliesToRange($logs.Range("D2:D27"))

$gFc = $logs.Range("G2:G26").FormatConditions
$gFc.Item(1).ModifyAppliesToRange($logs.Range("G2:G27"))
